# Update Name of Algo
# Apply updated KNN-imputed values in column B for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = 5.545
    14 = 5.755999999999999
    21 = 9.379000000000001
    23 = 7.398999999999999
    25 = 6.396
    26 = 6.547
    29 = 5.689
    53 = 6.032
    57 = 5.090999999999999
    59 = 4.435
    69 = 5.095
    79 = 5.623
    83 = 5.702
    91 = 5.468999999999999
    93 = 5.459
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
